# Trade #7 closed at 2026-02-17 04:06:34 - unknown UNKNOWN +0.000%
#
# A new (flat, break-even) closed trade for the MarketMaking strategy is
# appended as row 8 on both the "All Trades" and "MarketMaking" sheets, and
# the roll-up stats on "Summary" / "Strategy Status" are refreshed to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: Total P&L %, Total Trades, Win Rate %
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.09   # Total P&L %
$summary.Range("B6").Value = 7       # Total Trades
$summary.Range("B9").Value = 28.57   # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (Trades, Win Rate %)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 7        # Trades
$status.Range("G4").Value = 28.57    # Win Rate %

# ---------------------------------------------------------------------
# New trade row (Trade #7) appended to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 8

    $ws.Cells.Item($row, 1).Value = 7               # Trade #

    # Date/time look like dates to Excel's smart-parsing, so the written
    # text would otherwise get silently coerced into a date/time serial
    # number. Force them to stay as plain text (matching every other row
    # in the column), then drop back to the Normal style so no stray
    # number-format survives on the cell.
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"   # Date
    $ws.Cells.Item($row, 2).Style = "Normal"
    $ws.Cells.Item($row, 3).Value = "'04:06:28"     # Time
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"                      # Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"                              # Side
    $ws.Cells.Item($row, 6).Value = 0.79                                # Entry Price
    $ws.Cells.Item($row, 7).Value = 0.79                                # Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"                            # Status
    $ws.Cells.Item($row, 9).Value = 0                                   # P&L %
    $ws.Cells.Item($row, 10).Value = 0                                  # P&L $
    $ws.Cells.Item($row, 11).Value = 99.97                              # Capital After
    $ws.Cells.Item($row, 12).Value = 0                                  # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                                  # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                                # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps" # Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"                       # Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.12                               # Duration (min)
}
